# Generate Report for Handback
# Updates timestamps / status on the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-22 12:14:44"
$wsOverview.Range("G5").Value = "2016-08-22 12:14:44"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H4").Value = "2016-08-22 12:14:39"
$wsZhCn.Range("H5").Value = "2016-08-22 12:14:39"
$wsZhCn.Range("K4").Value = "2016-08-22 12:15:17"
$wsZhCn.Range("K5").Value = "2016-08-22 12:15:17"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H4").Value = "2016-08-22 12:14:44"
$wsDeDe.Range("H5").Value = "2016-08-22 12:14:44"
$wsDeDe.Range("K4").Value = "2016-08-22 12:15:23"
$wsDeDe.Range("K5").Value = "2016-08-22 12:15:23"
